$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$hyperlinkRunText = "linkedin.com/in/sivasankaran-pitchandi-273721182"
$oldSummaryPhrase = "Microsoft Certified: Data Analyst Associate and Power Platform Expert with nearly 3 years of experience in implementing analytical and reporting solutions. "
$newSummaryPhrase = "Microsoft Certified: Data Analyst Associate and Power Platform Expert with sound experience in implementing analytical and reporting solutions. "
$navyRgb = 6299648  # srgbClr val="002060" packed as a COM BGR int (0x00 60 20)

for ($shapeIdx = 1; $shapeIdx -le $s.Shapes.Count; $shapeIdx++) {
    $shp = $s.Shapes.Item($shapeIdx)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange

    # --- Contact-info textbox: color the "linkedin.com/..." hyperlink run and
    # --- the single space run right before it navy (srgbClr 002060). ---
    $runCount = $tr.Runs().Count
    for ($i = 1; $i -le $runCount; $i++) {
        $run = $tr.Runs($i)
        if ($run.Text -eq $hyperlinkRunText) {
            $run.Font.Color.RGB = $navyRgb
            if ($i -gt 1) {
                $prevRun = $tr.Runs($i - 1)
                $prevRun.Font.Color.RGB = $navyRgb
            }
        }
    }

    # --- Summary textbox: reword the certification blurb in place, keeping
    # --- the existing run formatting (single run, same rPr). ---
    $fullText = $tr.Text
    $pos = $fullText.IndexOf($oldSummaryPhrase)
    if ($pos -ge 0) {
        $target = $tr.Characters($pos + 1, $oldSummaryPhrase.Length)
        $target.Text = $newSummaryPhrase
    }
}
